$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 1.013382075887503
$ws.Range("C10").Value = 1.021144200595874
$ws.Range("C11").Value = 1.020265187051211
$ws.Range("C12").Value = 1.020136396917592
$ws.Range("C13").Value = 1.020421219193817
$ws.Range("C14").Value = 1.014783377783833
$ws.Range("C15").Value = 1.023577558129962
$ws.Range("C16").Value = 1.021533795909439
$ws.Range("C17").Value = 1.023089942607398
$ws.Range("C18").Value = 1.02318315857201
$ws.Range("C19").Value = 1.024902458094388
